# #5: insurance, claim, debt, investment done
#
# Fixes the "保險" (insurance) and "債務" (debt) sheets: the header row on
# both sheets was wrongly populated with data values instead of column
# names, and several columns present on the other sheets (property_category,
# category, date, legislator_name, legislator_id, source_file, index, ...)
# were missing entirely. This rebuilds both sheets with the full column set
# used across the rest of the workbook.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Sheet "保險" (insurance)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("保險")

# Extend the header formatting (bold + border, style of B1) and the plain
# data-row formatting (style of B2/B3) across the new columns before
# writing values into them.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("E1:K1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("E2:K3").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# Row 2 (index 81)
$ws.Range("B2").Value = "富邦人壽"
$ws.Range("C2").Value = "生存還本保險"
$ws.Range("D2").Value = "林正二"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
# "2012-04-12" looks like a date to Excel's auto-detection, so force it to
# be stored as text (matching the rest of the workbook) and then restore
# the plain "no explicit number format" styling from a sibling cell.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2012-04-12"
$ws.Range("F2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H2").Value = "林正二"
$ws.Range("I2").Value = 788
$ws.Range("J2").Value = "tmp32921"
$ws.Range("K2").Value = 81

# Row 3 (index 82)
$ws.Range("B3").Value = "富邦人壽"
$ws.Range("C3").Value = "年金保險"
$ws.Range("D3").Value = "林正二"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2012-04-12"
$ws.Range("F3").Copy() | Out-Null
$ws.Range("G3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H3").Value = "林正二"
$ws.Range("I3").Value = 788
$ws.Range("J3").Value = "tmp32921"
$ws.Range("K3").Value = 82

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet "債務" (debt)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("債務")

$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("H2:N4").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Row 2 (index 92)
$ws.Range("B2").Value = "房屋貸款"
$ws.Range("C2").Value = "林正二"
$ws.Range("D2").Value = "台灣土地銀行台東分行臺東縣台東市中華路"
$ws.Range("E2").Value = 430671
$ws.Range("F2").Value = "83年04月30日"
$ws.Range("G2").Value = "房屋貸款"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-12"
$ws.Range("I2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("K2").Value = "林正二"
$ws.Range("L2").Value = 788
$ws.Range("M2").Value = "tmp32921"
$ws.Range("N2").Value = 92

# Row 3 (index 93)
$ws.Range("B3").Value = "房屋貸款"
$ws.Range("C3").Value = "林正二"
$ws.Range("D3").Value = "台灣土地銀行花蓮分行花蓮縣花蓮市中山路"
$ws.Range("E3").Value = 4152495
$ws.Range("F3").Value = "88年04月07日"
$ws.Range("G3").Value = "房屋貸款"
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2012-04-12"
$ws.Range("I3").Copy() | Out-Null
$ws.Range("J3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("K3").Value = "林正二"
$ws.Range("L3").Value = 788
$ws.Range("M3").Value = "tmp32921"
$ws.Range("N3").Value = 93

# Row 4 (index 94)
$ws.Range("B4").Value = "房屋貸款"
$ws.Range("C4").Value = "林正二"
$ws.Range("D4").Value = "華南銀行花蓮分行花蓮縣花蓮市中山路"
$ws.Range("E4").Value = 2587777
$ws.Range("F4").Value = "88年04月01日"
$ws.Range("G4").Value = "房屋貸款"
$ws.Range("H4").Value = "debt"
$ws.Range("I4").Value = "normal"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "2012-04-12"
$ws.Range("I4").Copy() | Out-Null
$ws.Range("J4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("K4").Value = "林正二"
$ws.Range("L4").Value = 788
$ws.Range("M4").Value = "tmp32921"
$ws.Range("N4").Value = 94

$excel.CutCopyMode = $false
